$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric,
# so Excel stores them as text (matching the original inline-string type)
# instead of silently converting them to numbers.
$textForceCells = @("D5","D6","D7","D8","D9","D10","D12","D14","D15","D16","D19","D21","D22","D24","D26","D27","D29","D32","D33","D35","D36","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50","D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "37.137.79"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.022.31"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "246.31"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "59.66"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("D10").Value = "0.0806"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("D12").Value = "15.04"
$ws.Range("E12").Value = "  +6.38%  "
$ws.Range("D13").Value = "2.309.82"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "0.846"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "22.11"
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").Value = "5.43"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "2.015.78"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("D18").Value = "37.065.84"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "70.17"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "5.21"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").Value = "229.90"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "2.57"
$ws.Range("E24").Value = "  +5.21%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "9.34"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").Value = "164.31"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("D29").Value = "19.75"
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("E30").Value = "  +6.82%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "0.0668"
$ws.Range("D33").Value = "4.76"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  +13.31%  "
$ws.Range("D35").Value = "4.48"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "3.61"
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").Value = "0.0969"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").Value = "0.0216"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("D43").Value = "1.18"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "16.54"
$ws.Range("E44").Value = "  +4.29%  "
$ws.Range("D45").Value = "91.46"
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "7.54"
$ws.Range("E46").Value = "  +5.02%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.374.65"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "1.05"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  +15.32%  "
$ws.Range("D50").Value = "2.88"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").Value = "46.15"
$ws.Range("E51").Value = "  +0.13%  "

# Restore normal style on the text-forced cells (drop the temporary
# number-format override so styling matches the original workbook).
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
